$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite rows 4-36 with the updated localization entries (className.*/audioItem.*
# rows removed; disney.tips/disney.tips_desc and globalSearch.exercise inserted;
# everything below shifts up accordingly).
$ws.Range("A4").Value = 'disney.authorize_video'
$ws.Range("B4").Value = '英语视频'
$ws.Range("C4").Value = 'English video'
$ws.Range("D4").Value = '英語視頻'
$ws.Range("E4").Value = 'فيديو عربي'

$ws.Range("A5").Value = 'disney.skin_center'
$ws.Range("B5").Value = '皮肤中心'
$ws.Range("C5").Value = 'Skin Center'
$ws.Range("D5").Value = '皮膚中心'
$ws.Range("E5").Value = 'مركز الجلد'

$ws.Range("A6").Value = 'disney.more_theme'
$ws.Range("B6").Value = '更多主题'
$ws.Range("C6").Value = 'More'
$ws.Range("D6").Value = '更多主題'
$ws.Range("E6").Value = 'المزيد من المواضيع'

$ws.Range("A7").Value = 'disney.ipzone_placeholder'
$ws.Range("B7").Value = '搜索你想要的资源'
$ws.Range("C7").Value = 'Search for resources you want'
$ws.Range("D7").Value = '搜索你想要的資源'
$ws.Range("E7").Value = 'البحث عن الموارد التي تريد'

$ws.Range("A8").Value = 'disney.tips'
$ws.Range("B8").Value = '温馨提示'
$ws.Range("C8").Value = 'tips'
$ws.Range("D8").Value = '溫馨提示'
$ws.Range("E8").Value = 'نصائح'

$ws.Range("A9").Value = 'disney.tips_desc'
$ws.Range("B9").Value = '您的版本过旧，无法插入最新资源，快去更新后体验吧~'
$ws.Range("C9").Value = 'Your version is too old to insert the latest resources. Go to update and experience it~'
$ws.Range("D9").Value = '您的版本過舊，無法插入最新資源，快去更新後體驗吧~'
$ws.Range("E9").Value = 'أنت تستخدم إصدارًا منخفضًا جدًا. يُرجى ترقية الإصدار قبل استخدامه!'

$ws.Range("A10").Value = 'disney.update_now'
$ws.Range("B10").Value = '立即更新'
$ws.Range("C10").Value = 'Update now'
$ws.Range("D10").Value = '立即更新'
$ws.Range("E10").Value = 'تحديث فوري'

$ws.Range("A11").Value = 'disney.encourage_desc'
$ws.Range("B11").Value = '您可以在：放映态下 > 互动工具 > 表扬鼓励 > 表扬工具 中进行体验'
$ws.Range("C11").Value = 'You can experience it in: Show > Interactive Tools > praise and encouragement > praise tool'
$ws.Range("D11").Value = '您可以在：放映態下 > 互動工具 > 表揚鼓勵 > 表揚工具 中進行體驗'
$ws.Range("E11").Value = 'يمكنك تجربة مع أدوات تفاعلية في حالة العرض'

$ws.Range("A12").Value = 'disney.team_competition'
$ws.Range("B12").Value = '团队竞赛'
$ws.Range("C12").Value = 'Team Competition'
$ws.Range("D12").Value = '團隊競賽'
$ws.Range("E12").Value = 'فريق المنافسة'

$ws.Range("A13").Value = 'disney.random_rollcall'
$ws.Range("B13").Value = '随机点名'
$ws.Range("C13").Value = 'Random Roll Call'
$ws.Range("D13").Value = '隨機點名'
$ws.Range("E13").Value = 'اسم عشوائي'

$ws.Range("A14").Value = 'disney.class_encourage'
$ws.Range("B14").Value = '课堂鼓励'
$ws.Range("C14").Value = 'Classroom Encouragement'
$ws.Range("D14").Value = '課堂鼓勵'
$ws.Range("E14").Value = 'تشجيع الفصول الدراسية'

$ws.Range("A15").Value = 'disney.disneyResponder'
$ws.Range("B15").Value = '抢答工具'
$ws.Range("C15").Value = 'Answer Tool'
$ws.Range("D15").Value = '搶答工具'
$ws.Range("E15").Value = 'أداة الرد'

$ws.Range("A16").Value = 'disney.encourage'
$ws.Range("B16").Value = '鼓励'
$ws.Range("C16").Value = 'Encouragement'
$ws.Range("D16").Value = '鼓勵'
$ws.Range("E16").Value = 'شجع'

$ws.Range("A17").Value = 'disney.responder'
$ws.Range("B17").Value = '抢答'
$ws.Range("C17").Value = 'Answer'
$ws.Range("D17").Value = '搶答'
$ws.Range("E17").Value = 'الرد السريع'

$ws.Range("A18").Value = 'disney.fl_remark'
$ws.Range("B18").Value = '复仇者联盟邀你组队竞赛啦~'
$ws.Range("C18").Value = 'Avengers invite you to team up~'
$ws.Range("D18").Value = '復仇者聯盟邀妳組隊競賽啦~'
$ws.Range("E18").Value = 'دوري المنتقمون يدعوك إلى تنظيم المسابقة'

$ws.Range("A19").Value = 'disney.ml_remark'
$ws.Range("B19").Value = '可汗大点兵~抽学生课堂互动！'
$ws.Range("C19").Value = 'Khan big spot soldier ~ smoke student classroom interaction!'
$ws.Range("D19").Value = '可汗大點兵~抽學生課堂互動！'
$ws.Range("E19").Value = 'ما هو نوع من التفاعل الصفي ؟'

$ws.Range("A20").Value = 'disney.bx_remark'
$ws.Range("B20").Value = '放映态下>互动工具>表扬'
$ws.Range("C20").Value = 'Show > Interactive Tools > praise'
$ws.Range("D20").Value = '放映態下>互動工具>表揚'
$ws.Range("E20").Value = 'أداة تفاعلية في حالة العرض'

$ws.Range("A21").Value = 'disney.mq_remark'
$ws.Range("B21").Value = '与赛车手米奇一起抢答吧~'
$ws.Range("C21").Value = 'Join race car driver Mickey in the contest~'
$ws.Range("D21").Value = '與賽車手米奇壹起搶答吧~'
$ws.Range("E21").Value = 'مع ميكي ، متسابق'

$ws.Range("A22").Value = 'disney.use_now'
$ws.Range("B22").Value = '立即使用'
$ws.Range("C22").Value = 'Use it now'
$ws.Range("D22").Value = '立即使用'
$ws.Range("E22").Value = 'استخدام فوري'

$ws.Range("A23").Value = 'disney.login_use'
$ws.Range("B23").Value = '登录后即可使用'
$ws.Range("C23").Value = 'Please log in'
$ws.Range("D23").Value = '登錄後即可使用'
$ws.Range("E23").Value = 'بعد تسجيل الدخول يمكنك استخدام'

$ws.Range("A24").Value = 'disney.resource'
$ws.Range("B24").Value = '资源类型'
$ws.Range("C24").Value = 'Resource Type'
$ws.Range("D24").Value = '資源類型'
$ws.Range("E24").Value = 'أنواع الموارد'

$ws.Range("A25").Value = 'disney.current'
$ws.Range("B25").Value = '当前'
$ws.Range("C25").Value = 'Current'
$ws.Range("D25").Value = '當前'
$ws.Range("E25").Value = 'حالي'

$ws.Range("A26").Value = 'disney.ipzone'
$ws.Range("B26").Value = '迪士尼专区'
$ws.Range("C26").Value = 'Disney Zone'
$ws.Range("D26").Value = '迪士尼專區'
$ws.Range("E26").Value = 'ديزني لاند'

$ws.Range("A27").Value = 'globalSearch.exercise'
$ws.Range("B27").Value = '习题/试卷'
$ws.Range("C27").Value = 'Exercises & Test Papers'
$ws.Range("D27").Value = '習題/試卷'
$ws.Range("E27").Value = 'رياضي'

$ws.Range("A28").Value = 'authorWorks.activity_homepage'
$ws.Range("B28").Value = '活动首页'
$ws.Range("C28").Value = 'home'
$ws.Range("D28").Value = '活動首頁'
$ws.Range("E28").Value = 'الصفحة الرئيسية النشطة'

$ws.Range("A29").Value = 'authorWorks.my_works'
$ws.Range("B29").Value = '我的参赛作品'
$ws.Range("C29").Value = 'My entries'
$ws.Range("D29").Value = '我的參賽作品'
$ws.Range("E29").Value = 'بلدي يعمل'

$ws.Range("A30").Value = 'authorWorks.my_selected_works'
$ws.Range("B30").Value = '我的入选作品集'
$ws.Range("C30").Value = 'Selected works'
$ws.Range("D30").Value = '我的入選作品集'
$ws.Range("E30").Value = 'بلدي مجموعة مختارة'

$ws.Range("A31").Value = 'authorWorks.selected_works'
$ws.Range("B31").Value = '入选作品'
$ws.Range("C31").Value = 'selected'
$ws.Range("D31").Value = '入選作品'
$ws.Range("E31").Value = 'أعمال مختارة'

$ws.Range("A32").Value = 'authorWorks.award_works'
$ws.Range("B32").Value = '获奖作品'
$ws.Range("C32").Value = 'prize'
$ws.Range("D32").Value = '獲獎作品'
$ws.Range("E32").Value = 'الحائز على جائزة العمل'

$ws.Range("A33").Value = 'authorWorks.praise'
$ws.Range("B33").Value = '收获赞'
$ws.Range("C33").Value = 'praise'
$ws.Range("D33").Value = '收穫贊'
$ws.Range("E33").Value = 'حصاد الحمد'

$ws.Range("A34").Value = 'authorWorks.used'
$ws.Range("B34").Value = '被使用'
$ws.Range("C34").Value = 'used'
$ws.Range("D34").Value = '被使用'
$ws.Range("E34").Value = 'تستخدم'

$ws.Range("A35").Value = 'authorWorks.no_more'
$ws.Range("B35").Value = '没有更多了'
$ws.Range("C35").Value = 'no more'
$ws.Range("D35").Value = '沒有更多了'
$ws.Range("E35").Value = 'لا أكثر'

$ws.Range("A36").Value = 'authorWorks.take_part_in'
$ws.Range("B36").Value = '加入作品征集大赛，共建教学资源库'
$ws.Range("C36").Value = 'Join the competition of collecting works and build a teaching resource library'
$ws.Range("D36").Value = '加入作品徵集大賽，共建教學資源庫'
$ws.Range("E36").Value = 'الانضمام إلى مسابقة جمع الأعمال ، وبناء قاعدة بيانات الموارد التعليمية'

# Remove the now-unused trailing rows 37-42 (data shifted up by 6 rows; sheet
# shrinks from 42 to 36 rows).
$ws.Range("A37:E42").EntireRow.Delete()

Write-Host ("New dimension: " + $ws.UsedRange.Address())
